$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The committed change inserts one new data row at sheet row 239 ("Coliflor"
# price record dated 2022-09-05), pushing the former rows 239-324 down to
# 240-325 (dimension grows from A1:R324 to A1:R325).
$ws.Rows("239:239").Insert()

# Populate the newly inserted row 239 with its values.
$ws.Range("A239").Value = 5
$ws.Range("B239").Value = "Macroferia Regional de Talca"
$ws.Range("C239").Value = "Maule"
$ws.Range("D239").Value = 44809
$ws.Range("E239").Value = 7
$ws.Range("F239").Value = 100112008
$ws.Range("G239").Value = "Coliflor"
$ws.Range("H239").Value = "Sin especificar"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 2500
$ws.Range("K239").Value = 1300
$ws.Range("L239").Value = 1300
$ws.Range("M239").Value = 1300
$ws.Range("N239").Value = "$/unidad"
$ws.Range("O239").Value = "Región del Maule"
$ws.Range("P239").Value = 1300
$ws.Range("Q239").Value = 1
$ws.Range("R239").Value = "Hortaliza"
